$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Counts")

# --- Stage 1: capture every row-2 value we will need, before any cell gets
#              overwritten (several source/destination columns overlap).
# NB: use .Value2 for reads -- .Value (getter) does not resolve to the live
#     cell value in this host.
$old_B2 = $ws.Range("B2").Value2   # Non-Repository - on Course      -> new D2
$old_C2 = $ws.Range("C2").Value2   # Non-Repository Citation Matches -> new E2
$old_D2 = $ws.Range("D2").Value2   # Books on Course                 -> new F2
$old_E2 = $ws.Range("E2").Value2   # Physical Books on Course        -> new H2
$old_F2 = $ws.Range("F2").Value2   # No Electronic Version           -> new I2
$old_H2 = $ws.Range("H2").Value2   # Elec Already - Different Year   -> new J2
$old_I2 = $ws.Range("I2").Value2   # Elec Already - COVID            -> new K2
$old_J2 = $ws.Range("J2").Value2   # Elec Already - COVID Diff Year  -> new L2
$old_K2 = $ws.Range("K2").Value2   # Elec In Collection - Add        -> new M2
$old_L2 = $ws.Range("L2").Value2   # Elec In Collection - Potential  -> new N2
$old_M2 = $ws.Range("M2").Value2   # Elec Temporarily in Collection  -> new O2
$old_N2 = $ws.Range("N2").Value2   # Elec Temporarily - Diff Year    -> new P2
$old_O2 = $ws.Range("O2").Value2   # Course Code                     -> new B2
$old_P2 = $ws.Range("P2").Value2   # Course Name                     -> new C2
$old_Q2 = $ws.Range("Q2").Value2   # Electronic - Match on Course    -> new G2

# --- Stage 2: header row (row 1) -------------------------------------------------
# New layout inserts "Course Code"/"Course Name" right after "Processing
# Department", and inserts "Electronic - Match on Course" right after "Books on
# Course". The old "Electronic - Already on Course" (G1) and the two
# "Non-Match ..." columns (R1:S1) are dropped entirely.
$ws.Range("B1").Value = "Course Code"
$ws.Range("C1").Value = "Course Name"
$ws.Range("D1").Value = "Non-Repository - on Course"
$ws.Range("E1").Value = "Non-Repository Citation Matches"
$ws.Range("F1").Value = "Books on Course"
$ws.Range("G1").Value = "Electronic - Match on Course"
$ws.Range("H1").Value = "Physical Books on Course"
$ws.Range("I1").Value = "No Electronic Version for Physical Book"
$ws.Range("J1").Value = "Electronic - Already on Course - Different Year"
$ws.Range("K1").Value = "Electronic - Already on Course - COVID Temporary Electronic Collection"
$ws.Range("L1").Value = "Electronic - Already on Course - COVID Temporary Electronic Collection - Different Year"
$ws.Range("M1").Value = "Electronic - In Collection - Add to Course"
$ws.Range("N1").Value = "Electronic - In Collection - Potentially Add to Course - Different Year"
$ws.Range("O1").Value = "Electronic - Temporarily in Collection"
$ws.Range("P1").Value = "Electronic - Temporarily in Collection - Different Year"

# --- Stage 3: data row (row 2), written from the captured variables -------------
$ws.Range("B2").Value = $old_O2
$ws.Range("C2").Value = $old_P2
$ws.Range("D2").Value = $old_B2
$ws.Range("E2").Value = $old_C2
$ws.Range("F2").Value = $old_D2
$ws.Range("G2").Value = $old_Q2
$ws.Range("H2").Value = $old_E2
$ws.Range("I2").Value = $old_F2
$ws.Range("J2").Value = $old_H2
$ws.Range("K2").Value = $old_I2
$ws.Range("L2").Value = $old_J2
$ws.Range("M2").Value = $old_K2
$ws.Range("N2").Value = $old_L2
$ws.Range("O2").Value = $old_M2
$ws.Range("P2").Value = $old_N2

# --- Stage 4: drop the now-unused trailing columns Q:S (old "Electronic - Match
#              on Course" source, "Non-Match Ebooks on Course", "Non-Match
#              Temporary Collection Books on Course") --------------------------
$ws.Range("Q1:S1").EntireColumn.Delete()

# --- Stage 5: Totals row (row 4) - extend SUM formulas to the new K:P columns --
# K4/L4 already exist (they previously held the old M1:M2/N1:N2 sums) so simply
# overwriting their formula keeps their existing bold "Totals row" style (s="2").
# M4:P4 are brand-new cells in row 4, so they come in with the default (no)
# style, matching the diff.
$ws.Range("K4").Formula = "= SUM(K1:K2)"
$ws.Range("L4").Formula = "= SUM(L1:L2)"
$ws.Range("M4").Formula = "= SUM(M1:M2)"
$ws.Range("N4").Formula = "= SUM(N1:N2)"
$ws.Range("O4").Formula = "= SUM(O1:P2)"
$ws.Range("P4").Formula = "= SUM(P1:P2)"
